# Insert a new weekly price record for Coliflor at Macroferia Regional de
# Talca. The new record is inserted as row 68, pushing every existing
# record at/after row 68 down by one (row 154 -> row 155).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68 (shifts rows 68..154 down to 69..155).
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(68, 1).Value = 5
$ws.Cells.Item(68, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(68, 3).Value = "Maule"
$ws.Cells.Item(68, 4).Value = 44494
$ws.Cells.Item(68, 5).Value = 7
$ws.Cells.Item(68, 6).Value = 100112008
$ws.Cells.Item(68, 7).Value = "Coliflor"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 5000
$ws.Cells.Item(68, 11).Value = 600
$ws.Cells.Item(68, 12).Value = 600
$ws.Cells.Item(68, 13).Value = 600
$ws.Cells.Item(68, 14).Value = "$/unidad"
$ws.Cells.Item(68, 15).Value = "Región del Maule"
$ws.Cells.Item(68, 16).Value = 600
$ws.Cells.Item(68, 17).Value = 1
$ws.Cells.Item(68, 18).Value = "Hortaliza"
